# Update "想去人数" (column F) values for the rows that changed between
# data pulls, on both the "展览" and "全部类型" worksheets (which mirror
# the same rows).

$wb = $excel.ActiveWorkbook

$updates = @{
    "F6"  = 7356
    "F7"  = 466
    "F10" = 479
    "F11" = 14
    "F14" = 671
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
